$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 263
$ws.Range("I2").Value = 263
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 263
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -150
$ws.Range("N2").ClearContents()
# Row 6
$ws.Range("H6").Value = 905.0833
$ws.Range("J6").Value = 90
$ws.Range("L6").Value = 270
$ws.Range("N6").Value = -494
# Row 9
$ws.Range("H9").Value = 197.27272
$ws.Range("I9").Value = 212.44444
$ws.Range("K9").Value = 212.44444
$ws.Range("M9").Value = -43.44443999999999
# Row 10
$ws.Range("H10").Value = 225
$ws.Range("J10").Value = 225
$ws.Range("L10").Value = 225
$ws.Range("N10").Value = -811
# Row 11
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("K11").Value = 1
$ws.Range("M11").Value = 139
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
# Row 21
$ws.Range("H21").Value = 6135
$ws.Range("I21").Value = 6135
$ws.Range("K21").Value = 6135
$ws.Range("M21").Value = -5667
# Row 23
$ws.Range("H23").Value = 6135
$ws.Range("I23").Value = 6135
$ws.Range("K23").Value = 6135
$ws.Range("M23").Value = -5901
# Row 92
$ws.Range("H92").Value = 282.83334
$ws.Range("I92").Value = 282.83334
$ws.Range("K92").Value = 282.83334
$ws.Range("M92").Value = 965.16666
# Row 132
$ws.Range("H132").Value = 5514.0435
$ws.Range("I132").Value = 5514.0435
$ws.Range("K132").Value = 16542.1305
$ws.Range("M132").Value = -14012.1305
# Row 137
$ws.Range("H137").Value = 2000
$ws.Range("I137").Value = 2000
$ws.Range("K137").Value = 6000
$ws.Range("M137").Value = -3450
# Row 138
$ws.Range("H138").Value = 2654.875
$ws.Range("I138").Value = 754
$ws.Range("K138").Value = 2262
$ws.Range("M138").Value = 2878

$ws = $wb.Worksheets.Item("ARM")
# Row 18
$ws.Range("H18").Value = 5000
$ws.Range("J18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("N18").Value = -5644
# Row 61
$ws.Range("H61").Value = 4333.3335
$ws.Range("I61").Value = 4333.3335
$ws.Range("K61").Value = 4333.3335
$ws.Range("M61").Value = -4121.3335
# Row 136
$ws.Range("H136").Value = 4333.3335
$ws.Range("I136").Value = 4333.3335
$ws.Range("K136").Value = 13000.0005
$ws.Range("M136").Value = -10450.0005
# Row 139
$ws.Range("H139").Value = 60000
$ws.Range("I139").Value = 60000
$ws.Range("K139").Value = 60000
$ws.Range("M139").Value = -54860

$ws = $wb.Worksheets.Item("CRP")
# Row 69
$ws.Range("H69").Value = 41998.332
$ws.Range("I69").Value = 12997.5
$ws.Range("K69").Value = 12997.5
$ws.Range("M69").Value = -12248.5
# Row 72
$ws.Range("H72").Value = 41998.332
$ws.Range("I72").Value = 12997.5
$ws.Range("K72").Value = 38992.5
$ws.Range("M72").Value = -35248.5
# Row 105
$ws.Range("H105").Value = 1408.5
$ws.Range("I105").Value = 820.4286
$ws.Range("J105").Value = 2231.8
$ws.Range("K105").Value = 820.4286
$ws.Range("L105").Value = 2231.8
$ws.Range("M105").Value = 926.5714
$ws.Range("N105").Value = -5725.8
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 230
$ws.Range("I4").Value = 262.57144
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 787.71432
$ws.Range("L4").Value = 6
$ws.Range("M4").Value = -675.71432
$ws.Range("N4").Value = -230
# Row 15
$ws.Range("H15").Value = 50.8
$ws.Range("J15").Value = 50.8
$ws.Range("L15").Value = 152.4
$ws.Range("N15").Value = -432.4
# Row 41
$ws.Range("H41").Value = 1001
$ws.Range("I41").Value = 1001
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 3003
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -2665
$ws.Range("N41").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 28
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
# Row 97
$ws.Range("H97").Value = 428.33334
$ws.Range("I97").Value = 428.33334
$ws.Range("K97").Value = 428.33334
$ws.Range("M97").Value = 67.66665999999998

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 9051.200000000001
$ws.Range("I3").Value = 10114
$ws.Range("J3").Value = 4800
$ws.Range("K3").Value = 10114
$ws.Range("L3").Value = 4800
$ws.Range("M3").Value = -10002
$ws.Range("N3").Value = -5024
# Row 14
$ws.Range("H14").Value = 14005
$ws.Range("J14").Value = 14005
$ws.Range("L14").Value = 14005
$ws.Range("N14").Value = -14349
# Row 15
$ws.Range("H15").Value = 9051.200000000001
$ws.Range("I15").Value = 10114
$ws.Range("J15").Value = 4800
$ws.Range("K15").Value = 10114
$ws.Range("L15").Value = 4800
$ws.Range("M15").Value = -9944
$ws.Range("N15").Value = -5140
# Row 46
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
# Row 55
$ws.Range("H55").Value = 655.25
$ws.Range("I55").Value = 655.25
$ws.Range("K55").Value = 655.25
$ws.Range("M55").Value = -482.25
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 63
$ws.Range("H63").Value = 16900
$ws.Range("J63").Value = 16900
$ws.Range("L63").Value = 16900
$ws.Range("N63").Value = -18148
# Row 66
$ws.Range("H66").Value = 16900
$ws.Range("J66").Value = 16900
$ws.Range("L66").Value = 50700
$ws.Range("N66").Value = -56940
# Row 132
$ws.Range("H132").Value = 600
$ws.Range("I132").Value = 600
$ws.Range("K132").Value = 1800
$ws.Range("M132").Value = 730
